$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142477488380834
$ws.Range("D2").Value = 0.2412187059015309
$ws.Range("E2").Value = 0.1642394801235838
$ws.Range("F2").Value = 0.9298502608834838
$ws.Range("G2").Value = 0.4627726162836439
$ws.Range("H2").Value = 0.5607971153087448
$ws.Range("J2").Value = 0.1528597318075597
$ws.Range("M2").Value = 0.6290983407645285
$ws.Range("N2").Value = 1.549373067681017
$ws.Range("O2").Value = 2.002617514970126
$ws.Range("B3").Value = 0.1329197786953671
$ws.Range("D3").Value = 0.2428608296157826
$ws.Range("E3").Value = 0.1663019005649122
$ws.Range("F3").Value = 0.9175535045836938
$ws.Range("G3").Value = 0.4465487433260762
$ws.Range("H3").Value = 0.5571081219333536
$ws.Range("J3").Value = 0.155460469515436
$ws.Range("M3").Value = 0.5624023235514954
$ws.Range("N3").Value = 1.490262922841652
$ws.Range("O3").Value = 1.960176423489713
$ws.Range("B4").Value = 0.1271170574123914
$ws.Range("D4").Value = 0.2439697654102346
$ws.Range("E4").Value = 0.1676574136480324
$ws.Range("F4").Value = 0.9105867110871984
$ws.Range("G4").Value = 0.4369124299199711
$ws.Range("H4").Value = 0.5551391235319016
$ws.Range("J4").Value = 0.1571572748259129
$ws.Range("M4").Value = 0.5213437908738854
$ws.Range("N4").Value = 1.454474915804127
$ws.Range("O4").Value = 1.935414403037356
$ws.Range("B5").Value = 0.124769158632688
$ws.Range("D5").Value = 0.2444469962207947
$ws.Range("E5").Value = 0.1682322057337871
$ws.Range("F5").Value = 0.9078943358713687
$ws.Range("G5").Value = 0.4330671917606423
$ws.Range("H5").Value = 0.5544112227279925
$ws.Range("J5").Value = 0.1578738298574756
$ws.Range("M5").Value = 0.5045867486813052
$ws.Range("N5").Value = 1.44002040343787
$ws.Range("O5").Value = 1.925649586535371
$ws.Range("B6").Value = 0.1243803105060834
$ws.Range("D6").Value = 0.2445277704435291
$ws.Range("E6").Value = 0.1683290025555122
$ws.Range("F6").Value = 0.9074561243617865
$ws.Range("G6").Value = 0.4324336195655576
$ws.Range("H6").Value = 0.5542948549545343
$ws.Range("J6").Value = 0.1579943277189955
$ws.Range("M6").Value = 0.5018027680710304
$ws.Range("N6").Value = 1.437628122213255
$ws.Range("O6").Value = 1.924047823620128
$ws.Range("B7").Value = 0.1270853246820138
$ws.Range("D7").Value = 0.2439760989265096
$ws.Range("E7").Value = 0.1676650747800874
$ws.Range("F7").Value = 0.9105498070855802
$ws.Range("G7").Value = 0.4368602413304075
$ws.Range("H7").Value = 0.5551290051865152
$ws.Range("J7").Value = 0.1571668369874288
$ws.Range("M7").Value = 0.5211179001866242
$ws.Range("N7").Value = 1.45427945038179
$ws.Range("O7").Value = 1.935281392027406
$ws.Range("B8").Value = 0.1391684742673078
$ws.Range("D8").Value = 0.2417640274806416
$ws.Range("E8").Value = 0.1649320847175773
$ws.Range("F8").Value = 0.9254891943024006
$ws.Range("G8").Value = 0.4571110378262517
$ws.Range("H8").Value = 0.5594637221571617
$ws.Range("J8").Value = 0.1537357006306692
$ws.Range("M8").Value = 0.606124636958711
$ws.Range("N8").Value = 1.528888184525044
$ws.Range("O8").Value = 1.987714481841323
$ws.Range("B9").Value = 0.1633766678113773
$ws.Range("D9").Value = 0.2382241711541511
$ws.Range("E9").Value = 0.1602809269633898
$ws.Range("F9").Value = 0.9594207496415095
$ws.Range("G9").Value = 0.4994121276849341
$ws.Range("H9").Value = 0.5703128052223576
$ws.Range("J9").Value = 0.1478018578164679
$ws.Range("M9").Value = 0.7719151335328434
$ws.Range("N9").Value = 1.67912223153931
$ws.Range("O9").Value = 2.100845148501691
$ws.Range("B10").Value = 0.1814657063261222
$ws.Range("D10").Value = 0.2361090125121947
$ws.Range("E10").Value = 0.1572962843150698
$ws.Range("F10").Value = 0.9871885606793711
$ws.Range("G10").Value = 0.5320855005604699
$ws.Range("H10").Value = 0.5797167347012646
$ws.Range("J10").Value = 0.143928752475281
$ws.Range("M10").Value = 0.89310023300132
$ws.Range("N10").Value = 1.791788728889287
$ws.Range("O10").Value = 2.1902848345716
$ws.Range("B11").Value = 0.1897588163867709
$ws.Range("D11").Value = 0.2352520232532527
$ws.Range("E11").Value = 0.1560325925843742
$ws.Range("F11").Value = 1.000440082938738
$ws.Range("G11").Value = 0.5472994562836675
$ws.Range("H11").Value = 0.58430640341345
$ws.Range("J11").Value = 0.1422728537745197
$ws.Range("M11").Value = 0.9480819266081255
$ws.Range("N11").Value = 1.843519746742885
$ws.Range("O11").Value = 2.232354944450776
$ws.Range("B12").Value = 0.192908246638595
$ws.Range("D12").Value = 0.2349426173075813
$ws.Range("E12").Value = 0.1555676049017976
$ws.Range("F12").Value = 1.005547364230182
$ws.Range("G12").Value = 0.5531112243174192
$ws.Range("H12").Value = 0.5860892210660182
$ws.Range("J12").Value = 0.1416610909704215
$ws.Range("M12").Value = 0.9688796698089561
$ws.Range("N12").Value = 1.863175743631615
$ws.Range("O12").Value = 2.248485205284283
$ws.Range("B13").Value = 0.1922295635921074
$ws.Range("D13").Value = 0.2350085812497014
$ws.Range("E13").Value = 0.1556671455454435
$ws.Range("F13").Value = 1.004443450431012
$ws.Range("G13").Value = 0.5518573045089994
$ws.Range("H13").Value = 0.5857032668266129
$ws.Range("J13").Value = 0.1417921642417017
$ws.Range("M13").Value = 0.9644015361584053
$ws.Range("N13").Value = 1.858939546847608
$ws.Range("O13").Value = 2.245002400215185
$ws.Range("B14").Value = 0.1900177427154688
$ws.Range("D14").Value = 0.2352262653319812
$ws.Range("E14").Value = 0.1559940661780761
$ws.Range("F14").Value = 1.000858473482481
$ws.Range("G14").Value = 0.5477765795406242
$ws.Range("H14").Value = 0.5844521789727821
$ws.Range("J14").Value = 0.1422222169694756
$ws.Range("M14").Value = 0.9497934323531325
$ws.Range("N14").Value = 1.84513553565742
$ws.Range("O14").Value = 2.233677996446119
$ws.Range("B15").Value = 0.1886641047657349
$ws.Range("D15").Value = 0.235361571499169
$ws.Range("E15").Value = 0.1561960791252162
$ws.Range("F15").Value = 0.9986741929828327
$ws.Range("G15").Value = 0.5452836081186518
$ws.Range("H15").Value = 0.5836916864044781
$ws.Range("J15").Value = 0.1424876293277215
$ws.Range("M15").Value = 0.9408425513695136
$ws.Range("N15").Value = 1.836688785762675
$ws.Range("O15").Value = 2.226767424887697
$ws.Range("B16").Value = 0.1809250058408765
$ws.Range("D16").Value = 0.2361671348388086
$ws.Range("E16").Value = 0.1573807637356115
$ws.Range("F16").Value = 0.9863350223372152
$ws.Range("G16").Value = 0.5310983018066651
$ws.Range("H16").Value = 0.5794230607787796
$ws.Range("J16").Value = 0.1440391063659945
$ws.Range("M16").Value = 0.8895039719848796
$ws.Range("N16").Value = 1.788417422977147
$ws.Range("O16").Value = 2.187563317683185
$ws.Range("B17").Value = 0.1761936220116382
$ws.Range("D17").Value = 0.236688259812567
$ws.Range("E17").Value = 0.1581316305812681
$ws.Range("F17").Value = 0.9789241523801451
$ws.Range("G17").Value = 0.5224859967657949
$ws.Range("H17").Value = 0.5768842336363065
$ws.Range("J17").Value = 0.1450180696015568
$ws.Range("M17").Value = 0.8579708536247921
$ws.Range("N17").Value = 1.758925466310046
$ws.Range("O17").Value = 2.163867421130476
$ws.Range("B18").Value = 0.1734783231632093
$ws.Range("D18").Value = 0.236997900375691
$ws.Range("E18").Value = 0.1585723592113837
$ws.Range("F18").Value = 0.9747199497155208
$ws.Range("G18").Value = 0.5175654251587787
$ws.Range("H18").Value = 0.5754533144230294
$ws.Range("J18").Value = 0.1455911173728257
$ws.Range("M18").Value = 0.8398201915304782
$ws.Range("N18").Value = 1.742007639551304
$ws.Range("O18").Value = 2.150368378234305
$ws.Range("B19").Value = 0.1725600190960535
$ws.Range("D19").Value = 0.2371044405901301
$ws.Range("E19").Value = 0.1587231016417903
$ws.Range("F19").Value = 0.9733064950452928
$ws.Range("G19").Value = 0.5159050671890384
$ws.Range("H19").Value = 0.5749738708464207
$ws.Range("J19").Value = 0.1457868530871718
$ws.Range("M19").Value = 0.833672395320022
$ws.Range("N19").Value = 1.736287379449806
$ws.Range("O19").Value = 2.145820189672776
$ws.Range("B20").Value = 0.1766966596157147
$ws.Range("D20").Value = 0.2366317603017976
$ws.Range("E20").Value = 0.1580507834311078
$ws.Range("F20").Value = 0.9797070141590041
$ws.Range("G20").Value = 0.5233993750554475
$ws.Range("H20").Value = 0.5771514589791309
$ws.Range("J20").Value = 0.1449128246740887
$ws.Range("M20").Value = 0.8613290321144689
$ws.Range("N20").Value = 1.762060276996635
$ws.Range("O20").Value = 2.166376411542103
$ws.Range("B21").Value = 0.1906671655138297
$ws.Range("D21").Value = 0.2351619160515455
$ws.Range("E21").Value = 0.1558976738781119
$ws.Range("F21").Value = 1.001909046065734
$ws.Range("G21").Value = 0.5489738136786144
$ws.Range("H21").Value = 0.584818437861955
$ws.Range("J21").Value = 0.1420954847076672
$ws.Range("M21").Value = 0.954084813734184
$ws.Range("N21").Value = 1.849188318905476
$ws.Range("O21").Value = 2.236998839334774
$ws.Range("B22").Value = 0.1998501017880301
$ws.Range("D22").Value = 0.2342893977247797
$ws.Range("E22").Value = 0.1545694530821784
$ws.Range("F22").Value = 1.016939425757073
$ws.Range("G22").Value = 0.5659830340347014
$ws.Range("H22").Value = 0.5900904264960616
$ws.Range("J22").Value = 0.1403433348106518
$ws.Range("M22").Value = 1.014573414894244
$ws.Range("N22").Value = 1.906518611691126
$ws.Range("O22").Value = 2.284316027732302
$ws.Range("B23").Value = 0.1949442799805183
$ws.Range("D23").Value = 0.2347470188653702
$ws.Range("E23").Value = 0.1552711171703827
$ws.Range("F23").Value = 1.008869813462624
$ws.Range("G23").Value = 0.5568778664734992
$ws.Range("H23").Value = 0.587252777023167
$ws.Range("J23").Value = 0.1412703176031904
$ws.Range("M23").Value = 0.9823021871278854
$ws.Range("N23").Value = 1.875885666033412
$ws.Range("O23").Value = 2.258955601521677
$ws.Range("B24").Value = 0.176469221231585
$ws.Range("D24").Value = 0.2366572724598655
$ws.Range("E24").Value = 0.1580873062541954
$ws.Range("F24").Value = 0.9793529066551656
$ws.Range("G24").Value = 0.5229863409344944
$ws.Range("H24").Value = 0.5770305570410414
$ws.Range("J24").Value = 0.1449603740471481
$ws.Range("M24").Value = 0.8598108675611087
$ws.Range("N24").Value = 1.760642912132823
$ws.Range("O24").Value = 2.165241710603055
$ws.Range("B25").Value = 0.1567737604999309
$ws.Range("D25").Value = 0.2390964583580555
$ws.Range("E25").Value = 0.1614632841919335
$ws.Range("F25").Value = 0.9497438226178048
$ws.Range("G25").Value = 0.4876896733852618
$ws.Range("H25").Value = 0.5671262642298842
$ws.Range("J25").Value = 0.1493218371077276
$ws.Range("M25").Value = 0.727168609484778
$ws.Range("N25").Value = 1.638070420300863
$ws.Range("O25").Value = 2.06913270280009
